$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = '< -0.0273'
$ws.Range('E2').Value = '0.0851 -0.0468'
$ws.Range('H2').Value = '< -0.063'
$ws.Range('J2').Value = '< -0.0762'
$ws.Range('L2').Value = '< -0.0762'
$ws.Range('D3').Value = '< -0.0489'
$ws.Range('E3').Value = '0.2749 -0.0194'
$ws.Range('F3').Value = '< -0.0489'
$ws.Range('G3').Value = '0.6566 0.014'
$ws.Range('J3').Value = '< -0.0489'
$ws.Range('K3').Value = '0.0265 -0.04'
$ws.Range('L3').Value = '< -0.0489'
$ws.Range('M3').Value = '< -0.0489'
$ws.Range('N3').Value = '0.0139 -0.04'
$ws.Range('O3').Value = '< -0.0489'
$ws.Range('E4').Value = '0.7038 0.0295'
$ws.Range('H4').Value = '0.7275 0.0132'
$ws.Range('I4').Value = '0.7125 0.0088'
$ws.Range('J4').Value = '0.5784 0'
$ws.Range('K4').Value = '0.6847 0.0089'
$ws.Range('L4').Value = '0.5784 0'
$ws.Range('N4').Value = '0.7982 0.0089'
$ws.Range('F5').Value = '0.295 -0.0295'
$ws.Range('G5').Value = '0.8622 0.0335'
$ws.Range('H5').Value = '0.5345 -0.0162'
$ws.Range('I5').Value = '0.5063 -0.0207'
$ws.Range('J5').Value = '0.4103 -0.0295'
$ws.Range('K5').Value = '0.4436 -0.0206'
$ws.Range('L5').Value = '0.4103 -0.0295'
$ws.Range('M5').Value = '0.295 -0.0295'
$ws.Range('N5').Value = '0.4933 -0.0206'
$ws.Range('O5').Value = '0.295 -0.0295'
$ws.Range('H6').Value = '0.7275 0.0132'
$ws.Range('I6').Value = '0.7125 0.0088'
$ws.Range('J6').Value = '0.5784 0'
$ws.Range('K6').Value = '0.6847 0.0089'
$ws.Range('L6').Value = '0.5784 0'
$ws.Range('N6').Value = '0.7982 0.0089'
$ws.Range('H7').Value = '< -0.0497'
$ws.Range('J8').Value = '0.0036 -0.0132'
$ws.Range('K8').Value = '0.3325 -0.0044'
$ws.Range('L8').Value = '0.0036 -0.0132'
$ws.Range('M8').Value = '0.2331 -0.0132'
$ws.Range('N8').Value = '0.3202 -0.0044'
$ws.Range('O8').Value = '0.2331 -0.0132'
$ws.Range('J9').Value = '0.0349 -0.0088'
$ws.Range('K9').Value = '0.4337 1e-04'
$ws.Range('L9').Value = '0.0349 -0.0088'
$ws.Range('M9').Value = '0.2481 -0.0088'
$ws.Range('N9').Value = '0.5141 1e-04'
$ws.Range('O9').Value = '0.2481 -0.0088'
$ws.Range('K10').Value = '0.5026 0.0089'
$ws.Range('M10').Value = '0.3822 0'
$ws.Range('N10').Value = '0.8125 0.0089'
$ws.Range('O10').Value = '0.3822 0'
$ws.Range('L11').Value = '0.4766 -0.0089'
$ws.Range('M11').Value = '0.2796 -0.0089'
$ws.Range('N11').Value = '0.6927 0'
$ws.Range('O11').Value = '0.2796 -0.0089'
$ws.Range('M12').Value = '0.3822 0'
$ws.Range('N12').Value = '0.8125 0.0089'
$ws.Range('O12').Value = '0.3822 0'
$ws.Range('N13').Value = '0.7982 0.0089'
$ws.Range('O14').Value = '0.181 -0.0089'
